$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "activo" column values from "TRUE" to "Activo"
$ws.Range("D2").Value = "Activo"
$ws.Range("D3").Value = "Activo"
$ws.Range("D4").Value = "Activo"
$ws.Range("D5").Value = "Activo"

# Reflect the user's selection over the edited range
$ws.Range("D2:D5").Select() | Out-Null
